$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.967.62"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +16.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.665.26"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +13.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9969"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.93"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +11.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9902"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3734"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +5.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3466"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +12.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.33"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +22.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.176"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +8.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07253"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +9.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9930"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.80"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +14.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.035"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +10.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.746"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +9.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.664.80"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +13.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001103"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +8.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9897"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06729"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +12.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "81.39"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +18.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.48"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +13.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.080"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +11.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.01"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +7.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.955.13"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +16.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.359"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.51%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.696"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +29.15%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.382"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -9.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.10"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.36%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +13.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.840.57"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +13.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.07"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +11.18%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.247"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +26.68%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.068"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9861"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +24.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.722"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +19.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08393"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "12.37"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +20.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.986"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +23.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06367"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +12.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.317"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +12.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.288"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02324"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +14.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2075"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +12.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6097"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +16.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9889"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.816"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +8.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.26"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +10.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5952"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +15.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "127.22"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.005"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +11.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07084"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +10.59%  "
